# fix(FN-3460): fix invalid facility utilisation values -- all need to match as same facility id for all rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 - correct utilisation values so facility id matches across all rows
$ws.Range("E5").Value = 600000
$ws.Range("F5").Value = 761579.37
$ws.Range("G5").Value = 3938753.8
# H5 left unchanged (456)

# Row 6 - correct utilisation values so facility id matches across all rows
$ws.Range("E6").Value = 600000
# F6 left unchanged (761579.37)
$ws.Range("G6").Value = 761579.37
# H6 left unchanged (456.77)

# Adjust column widths: columns E:G (5-7) now share the same width/bestFit formatting
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 16.33203125
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 16.33203125

# Update the active selection to reflect the edited range
$ws.Range("E5:H6").Select()
